$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 2695.1667
$ws.Range("I19").Value = 2435.4285
$ws.Range("J19").Value = 2860.4546
$ws.Range("K19").Value = 2435.4285
$ws.Range("L19").Value = 2860.4546
$ws.Range("M19").Value = -2260.4285
$ws.Range("N19").Value = -3210.4546
$ws.Range("H28").Value = 39423.445
$ws.Range("I28").Value = 48356.617
$ws.Range("K28").Value = 48356.617
$ws.Range("M28").Value = -47871.617
$ws.Range("H58").Value = 1878.4286
$ws.Range("I58").Value = 177.55556
$ws.Range("J58").Value = 4940
$ws.Range("K58").Value = 532.66668
$ws.Range("L58").Value = 14820
$ws.Range("M58").Value = -382.66668
$ws.Range("N58").Value = -15120
$ws.Range("H135").Value = 111117896
$ws.Range("I135").Value = 111117896
$ws.Range("K135").Value = 1000061064
$ws.Range("M135").Value = -1000058529
$ws.Range("H137").Value = 827.4375
$ws.Range("I137").Value = 827.4375
$ws.Range("K137").Value = 2482.3125
$ws.Range("M137").Value = 67.6875
$ws.Range("H138").Value = 2516.0286
$ws.Range("I138").Value = 1459.8422
$ws.Range("K138").Value = 4379.5266
$ws.Range("M138").Value = 760.4733999999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1252.5555
$ws.Range("I2").Value = 1457.1428
$ws.Range("K2").Value = 1457.1428
$ws.Range("M2").Value = -1344.1428
$ws.Range("H4").Value = 7519.2856
$ws.Range("I4").Value = 278.2857
$ws.Range("K4").Value = 278.2857
$ws.Range("M4").Value = -162.2857
$ws.Range("H10").Value = 955
$ws.Range("J10").Value = 955
$ws.Range("L10").Value = 955
$ws.Range("N10").Value = -1295
$ws.Range("H32").Value = 2255.4219
$ws.Range("I32").Value = 1369.585
$ws.Range("K32").Value = 1369.585
$ws.Range("M32").Value = -1082.585
$ws.Range("H45").Value = 2068.5
$ws.Range("I45").Value = 1627.75
$ws.Range("K45").Value = 1627.75
$ws.Range("M45").Value = -1250.75
$ws.Range("H61").Value = 30304310
$ws.Range("I61").Value = 37038156
$ws.Range("K61").Value = 37038156
$ws.Range("M61").Value = -37037944
$ws.Range("H74").Value = 33336322
$ws.Range("I74").Value = 45456532
$ws.Range("J74").Value = 5744.5
$ws.Range("K74").Value = 45456532
$ws.Range("L74").Value = 5744.5
$ws.Range("M74").Value = -45455658
$ws.Range("N74").Value = -7492.5
$ws.Range("H77").Value = 33336322
$ws.Range("I77").Value = 45456532
$ws.Range("J77").Value = 5744.5
$ws.Range("K77").Value = 227282660
$ws.Range("L77").Value = 28722.5
$ws.Range("M77").Value = -227278292
$ws.Range("N77").Value = -37458.5
$ws.Range("H102").Value = 2032639.9
$ws.Range("I102").Value = 2165834.5
$ws.Range("J102").Value = 167916.33
$ws.Range("K102").Value = 2165834.5
$ws.Range("L102").Value = 167916.33
$ws.Range("M102").Value = -2164212.5
$ws.Range("N102").Value = -171160.33
$ws.Range("H110").Value = 125005050
$ws.Range("I110").Value = 200005580
$ws.Range("J110").Value = 4140
$ws.Range("K110").Value = 200005580
$ws.Range("L110").Value = 4140
$ws.Range("M110").Value = -200003535
$ws.Range("N110").Value = -8230
$ws.Range("H116").Value = 1252.5555
$ws.Range("I116").Value = 1457.1428
$ws.Range("K116").Value = 1457.1428
$ws.Range("M116").Value = 836.8571999999999
$ws.Range("H122").Value = 22225238
$ws.Range("I122").Value = 37038840
$ws.Range("K122").Value = 111116520
$ws.Range("M122").Value = -111114070
$ws.Range("H132").Value = 22729072
$ws.Range("I132").Value = 25642704
$ws.Range("J132").Value = 2740
$ws.Range("K132").Value = 76928112
$ws.Range("L132").Value = 8220
$ws.Range("M132").Value = -76925582
$ws.Range("N132").Value = -13280
$ws.Range("H136").Value = 30304310
$ws.Range("I136").Value = 37038156
$ws.Range("K136").Value = 111114468
$ws.Range("M136").Value = -111111918

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1252.5555
$ws.Range("I3").Value = 1457.1428
$ws.Range("K3").Value = 1457.1428
$ws.Range("M3").Value = -1343.1428
$ws.Range("H105").Value = 1922.25
$ws.Range("I105").Value = 1843.6471
$ws.Range("J105").Value = 2113.1428
$ws.Range("K105").Value = 1843.6471
$ws.Range("L105").Value = 2113.1428
$ws.Range("M105").Value = -96.64709999999991
$ws.Range("N105").Value = -5607.1428
$ws.Range("H134").Value = 2174.8235
$ws.Range("I134").Value = 1165.6666
$ws.Range("J134").Value = 9743.5
$ws.Range("K134").Value = 3496.9998
$ws.Range("L134").Value = 29230.5
$ws.Range("M134").Value = -961.9998000000001
$ws.Range("N134").Value = -34300.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H5").Value = 1029.875
$ws.Range("I5").Value = 706.5
$ws.Range("J5").Value = 2000
$ws.Range("K5").Value = 706.5
$ws.Range("L5").Value = 2000
$ws.Range("M5").Value = -594.5
$ws.Range("N5").Value = -2224
$ws.Range("H7").Value = 89.833336
$ws.Range("I7").Value = 97.86667
$ws.Range("K7").Value = 97.86667
$ws.Range("M7").Value = 15.13333
$ws.Range("H107").Value = 3165.8572
$ws.Range("I107").Value = 1774.7273
$ws.Range("K107").Value = 1774.7273
$ws.Range("M107").Value = 145.2727
$ws.Range("H124").Value = 58000
$ws.Range("J124").Value = 58000
$ws.Range("L124").Value = 58000
$ws.Range("N124").Value = -62910
$ws.Range("H132").Value = 3064.889
$ws.Range("I132").Value = 2916.8125
$ws.Range("K132").Value = 8750.4375
$ws.Range("M132").Value = -6220.4375
$ws.Range("H134").Value = 2020.3334
$ws.Range("I134").Value = 1528.3334
$ws.Range("J134").Value = 2319.0715
$ws.Range("K134").Value = 4585.0002
$ws.Range("L134").Value = 9012.999899999999
$ws.Range("M134").Value = -2050.0002

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 4516.3335
$ws.Range("I56").Value = 4516.3335
$ws.Range("K56").Value = 4516.3335
$ws.Range("M56").Value = -3986.3335
$ws.Range("H75").Value = 781.8182
$ws.Range("J75").Value = 1055.4286
$ws.Range("L75").Value = 3166.2858
$ws.Range("N75").Value = -5162.2858
$ws.Range("H78").Value = 781.8182
$ws.Range("J78").Value = 1055.4286
$ws.Range("L78").Value = 9498.857399999999
$ws.Range("N78").Value = -19482.8574
$ws.Range("H131").Value = 38176.168
$ws.Range("J131").Value = 45631.4
$ws.Range("L131").Value = 136894.2
$ws.Range("N131").Value = -146974.2
$ws.Range("H137").Value = 41793350
$ws.Range("I137").Value = 41793350
$ws.Range("K137").Value = 125380050
$ws.Range("M137").Value = -125374950

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H40").Value = 3432.6
$ws.Range("I40").Value = 1500
$ws.Range("J40").Value = 3915.75
$ws.Range("K40").Value = 1500
$ws.Range("L40").Value = 3915.75
$ws.Range("M40").Value = -1349
$ws.Range("N40").Value = -4217.75
$ws.Range("H102").Value = 2111.9688
$ws.Range("I102").Value = 1317.8422
$ws.Range("K102").Value = 1317.8422
$ws.Range("M102").Value = 304.1578
$ws.Range("H107").Value = 2293.125
$ws.Range("J107").Value = 2490.5
$ws.Range("L107").Value = 2490.5
$ws.Range("N107").Value = -6330.5
$ws.Range("H113").Value = 2129.2
$ws.Range("I113").Value = 1887.5454
$ws.Range("J113").Value = 2319.0715
$ws.Range("K113").Value = 1887.5454
$ws.Range("L113").Value = 2319.0715
$ws.Range("M113").Value = 282.4546
$ws.Range("N113").Value = -6659.0715
$ws.Range("H122").Value = 2521.25
$ws.Range("I122").Value = 2286.6086
$ws.Range("K122").Value = 6859.825800000001
$ws.Range("M122").Value = -4409.825800000001
$ws.Range("H132").Value = 3161.1738
$ws.Range("I132").Value = 3006.5151
$ws.Range("J132").Value = 3553.7693
$ws.Range("K132").Value = 9019.5453
$ws.Range("L132").Value = 10661.3079
$ws.Range("M132").Value = -6489.5453
$ws.Range("N132").Value = -15721.3079

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1000
$ws.Range("J22").Value = 1000
$ws.Range("L22").Value = 1000
$ws.Range("N22").Value = -1590
$ws.Range("H27").Value = 1000
$ws.Range("J27").Value = 1000
$ws.Range("L27").Value = 1000
$ws.Range("N27").Value = -1214
$ws.Range("H55").Value = 494.4375
$ws.Range("I55").Value = 474.8
$ws.Range("J55").Value = 503.36365
$ws.Range("K55").Value = 474.8
$ws.Range("L55").Value = 503.36365
$ws.Range("M55").Value = -301.8
$ws.Range("N55").Value = -849.36365
$ws.Range("H82").Value = 3247.75
$ws.Range("I82").Value = 3247.75
$ws.Range("K82").Value = 3247.75
$ws.Range("M82").Value = -2886.75
$ws.Range("H85").Value = 3247.75
$ws.Range("I85").Value = 3247.75
$ws.Range("K85").Value = 3247.75
$ws.Range("M85").Value = -1999.75
$ws.Range("H132").Value = 6529.1904
$ws.Range("I132").Value = 3439.3076
$ws.Range("J132").Value = 11550.25
$ws.Range("K132").Value = 10317.9228
$ws.Range("L132").Value = 34650.75
$ws.Range("M132").Value = -7787.9228
$ws.Range("N132").Value = -39710.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H58").Value = 10000
$ws.Range("J58").Value = 10000
$ws.Range("L58").Value = 10000
$ws.Range("N58").Value = -10616
$ws.Range("H107").Value = 571.9722
$ws.Range("I107").Value = 460.6087
$ws.Range("J107").Value = 769
$ws.Range("K107").Value = 1381.8261
$ws.Range("L107").Value = 2307
$ws.Range("M107").Value = 538.1739
$ws.Range("N107").Value = -6147
$ws.Range("H132").Value = 4929.077
$ws.Range("I132").Value = 4712
$ws.Range("J132").Value = 5652.6665
$ws.Range("K132").Value = 14136
$ws.Range("L132").Value = 16957.9995
$ws.Range("M132").Value = -11606
$ws.Range("N132").Value = -22017.9995
$ws.Range("H136").Value = 5908.273
$ws.Range("I136").Value = 3497.5
$ws.Range("J136").Value = 6444
$ws.Range("K136").Value = 10492.5
$ws.Range("L136").Value = 19332
$ws.Range("M136").Value = -7942.5
$ws.Range("N136").Value = -24432
